$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B11").Value = 9
$ws.Range("C11").Value = 2
$ws.Range("B12").Value = 171
$ws.Range("C12").Value = -12
$ws.Range("E12").Value = "159/252"
